$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header cell A1 previously read "currencies[0].id" (an array-style JSON
# path for a single-currency schema). Rename it to the flat "currencies.id"
# form used by the new base JSON schema (node-xlsx).
$ws.Range("A1").Value = "currencies.id"

# Reflect the author's final selection/scroll state: cell B1 selected and
# the view scrolled back so column A is visible (no more forced topLeftCell).
$ws.Range("B1").Select()
